$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new precondition sentence to the existing precondition text (B9).
$ws.Range("B9").Value = "Que el actor tenga los permisos necesarios para ver el registro." + [char]10 + "Que existan oportunidades asignadas a telemarketers."

# The cell now wraps to two lines; match the row height Excel would compute for it.
$ws.Rows.Item(9).RowHeight = 25.5

# Move the active selection to C1 (matches the saved view state in the edit).
$ws.Range("C1").Select() | Out-Null
